$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / safe-string assignments (values that Excel will not
# auto-convert to numbers, e.g. contain letters, %, multiple dots, or spaces).
$plainUpdates = @(
    @('D2', '67.670.26'),
    @('E2', '  -1.04%  '),
    @('D3', '3.786.30'),
    @('E3', '  -0.14%  '),
    @('E4', '  +0.10%  '),
    @('E5', '  +0.11%  '),
    @('E6', '  -0.82%  '),
    @('D7', '3.785.61'),
    @('E7', '  -0.12%  '),
    @('E8', '  +0.02%  '),
    @('E9', '  -0.18%  '),
    @('E10', '  -0.94%  '),
    @('E11', '  -1.00%  '),
    @('E12', '  -0.39%  '),
    @('E13', '  -2.69%  '),
    @('E14', '  -0.53%  '),
    @('D15', '4.423.72'),
    @('E15', '  +0.05%  '),
    @('D16', '3.797.36'),
    @('E16', '  +0.12%  '),
    @('E17', '  +3.13%  '),
    @('D18', '67.644.98'),
    @('E18', '  -1.04%  '),
    @('E19', '  +1.28%  '),
    @('E20', '  +0.10%  '),
    @('E21', '  -9.11%  '),
    @('E22', '  -1.49%  '),
    @('E23', '  -0.24%  '),
    @('E24', '  +2.47%  '),
    @('E26', '  +0.89%  '),
    @('E27', '  -3.12%  '),
    @('E28', '  -0.04%  '),
    @('E29', '  -1.79%  '),
    @('D30', '3.935.96'),
    @('E30', '  +0.02%  '),
    @('E31', '  -0.04%  '),
    @('E32', '  +4.45%  '),
    @('E33', '  -1.42%  '),
    @('E34', '  -1.77%  '),
    @('B35', 'Binance-PegBSC-USD'),
    @('C35', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'),
    @('E35', '  -0.39%  '),
    @('B36', 'Aptos'),
    @('C36', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'),
    @('E36', '  -1.15%  '),
    @('E37', '  -0.84%  '),
    @('E38', '  -2.57%  '),
    @('E39', '  -0.46%  '),
    @('E40', '  -1.06%  '),
    @('E41', '  -0.56%  '),
    @('E45', '  -0.29%  '),
    @('E46', '  -1.58%  '),
    @('E47', '  +2.79%  '),
    @('E48', '  -1.78%  '),
    @('E49', '  +5.34%  '),
    @('E50', '  -1.54%  '),
    @('E51', '  -5.00%  ')
)

foreach ($pair in $plainUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Numeric-looking strings that must stay text (column D prices without
# thousands separators collide with Excel general number parsing).
# Force them to remain text by writing a text formula, then collapsing
# the formula down to a static value via Copy + PasteSpecial(xlPasteValues),
# which does not touch the cell style (unlike NumberFormat = "@").
$forceTextUpdates = @(
    @('D5', '595.70'),
    @('D6', '166.71'),
    @('D11', '6.35'),
    @('D13', '0.0000252'),
    @('D14', '35.99'),
    @('D17', '18.54'),
    @('D22', '458.71'),
    @('D23', '0.698'),
    @('D24', '0.0000152'),
    @('D25', '83.38'),
    @('D26', '12.04'),
    @('D27', '2.10'),
    @('D32', '2.26'),
    @('D33', '7.21'),
    @('D34', '29.54'),
    @('D35', '0.998'),
    @('D36', '9.05'),
    @('D37', '0.0999'),
    @('D38', '3.34'),
    @('D41', '5.76'),
    @('D44', '48.02'),
    @('D45', '43.85'),
    @('D46', '0.296'),
    @('D47', '150.72'),
    @('D48', '8.27'),
    @('D49', '26.76'),
    @('D50', '388.56')
)

foreach ($pair in $forceTextUpdates) {
    $ref = $pair[0]
    $val = $pair[1]
    $cell = $ws.Range($ref)
    $cell.Formula = '="' + $val + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
